$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure the updated cells keep their original "number stored as text" type
# (matching the source t="str" cells) rather than being auto-converted to numbers.
$dataRange = $ws.Range("C2:F13")
$dataRange.NumberFormat = "@"

$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "1"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "0"

$ws.Range("C3").Value = "13"
$ws.Range("D3").Value = "10"
$ws.Range("E3").Value = "0"
$ws.Range("F3").Value = "1"

$ws.Range("C4").Value = "34"
$ws.Range("D4").Value = "30"
$ws.Range("E4").Value = "4"
$ws.Range("F4").Value = "1"

$ws.Range("C5").Value = "0"
$ws.Range("D5").Value = "0"
$ws.Range("E5").Value = "0"
$ws.Range("F5").Value = "0"

$ws.Range("C6").Value = "0"
$ws.Range("D6").Value = "3"
$ws.Range("E6").Value = "0"
$ws.Range("F6").Value = "0"

$ws.Range("C7").Value = "12"
$ws.Range("D7").Value = "7"
$ws.Range("E7").Value = "2"
$ws.Range("F7").Value = "0"

$ws.Range("C8").Value = "1"
$ws.Range("D8").Value = "3"
$ws.Range("E8").Value = "0"
$ws.Range("F8").Value = "0"

$ws.Range("C9").Value = "3"
$ws.Range("D9").Value = "4"
$ws.Range("E9").Value = "0"
$ws.Range("F9").Value = "0"

$ws.Range("C10").Value = "3"
$ws.Range("D10").Value = "3"
$ws.Range("E10").Value = "0"
$ws.Range("F10").Value = "0"

$ws.Range("C11").Value = "12"
$ws.Range("D11").Value = "17"
$ws.Range("E11").Value = "0"
$ws.Range("F11").Value = "1"

$ws.Range("C12").Value = "10"
$ws.Range("D12").Value = "10"
$ws.Range("E12").Value = "1"
$ws.Range("F12").Value = "0"

$ws.Range("C13").Value = "20"
$ws.Range("D13").Value = "4"
$ws.Range("E13").Value = "2"
$ws.Range("F13").Value = "2"

